$wb = $excel.ActiveWorkbook

# Avoid confirmation prompt when deleting a worksheet
$excel.DisplayAlerts = $false

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete the "Desarquivamentos Pendentes" sheet entirely
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true
